$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (row 1), matching the formatting (bold,
# centered, bordered) already used by the other header cells such as L1.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Populate the three new data columns for every existing data row (2-25).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"
    $ws.Cells.Item($r, 14).Value = 20181295
    $ws.Cells.Item($r, 15).Value = 7
}
